# edit.ps1 - applies the cover_sheet_2.0.docx template fix:
#  1. Picture "Graphic 5" alt text: "A cell phone..." -> "Mobile phone..."
#     (also refreshes the drawing's editId and the inner pic:cNvPr id, matching
#     what Word does when the alt text is edited through the UI)
#  2. Paragraph spacing/formatting fix around the "{%p endif %}" template tag:
#     the tag now lives in the paragraph that used to be blank (with the
#     mismatched rFonts/szCs), and the paragraph that used to hold the tag is
#     now blank with plain spacing.

$d = $word.ActiveDocument

function Get-ParagraphIndexAt($doc, $pos) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($pos -ge $p.Range.Start -and $pos -lt $p.Range.End) {
            return $i
        }
    }
    return -1
}

# --- Part 1: fix the picture's alt text (and the ids Word keeps in sync with it) ---

$anchorXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="0CEC689E" w14:textId="10007574" w:rsidR="00CC26EE" w:rsidRDefault="00671B5C" w:rsidP="00064D95"><w:pPr><w:spacing w:after="0"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00B42820"><w:rPr><w:noProof/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251672576" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="0" wp14:anchorId="7E59310E" wp14:editId="534FB8A8"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="margin"><wp:posOffset>4977130</wp:posOffset></wp:positionH><wp:positionV relativeFrom="page"><wp:posOffset>3307080</wp:posOffset></wp:positionV><wp:extent cx="987552" cy="1216152"/><wp:effectExtent l="0" t="0" r="3175" b="3175"/><wp:wrapThrough wrapText="bothSides"><wp:wrapPolygon edited="0"><wp:start x="0" y="0"/><wp:lineTo x="0" y="21318"/><wp:lineTo x="14585" y="21318"/><wp:lineTo x="15002" y="21318"/><wp:lineTo x="15835" y="16242"/><wp:lineTo x="21253" y="12182"/><wp:lineTo x="21253" y="3722"/><wp:lineTo x="15002" y="0"/><wp:lineTo x="0" y="0"/></wp:wrapPolygon></wp:wrapThrough><wp:docPr id="258308898" name="Graphic 5" descr="Mobile phone with the words What's Next Text in a speech bubble."><a:extLst xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:ext uri="{C183D7F6-B498-43B3-948B-1728B52AA6E4}"><adec:decorative xmlns:adec="http://schemas.microsoft.com/office/drawing/2017/decorative" val="0"/></a:ext></a:extLst></wp:docPr><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="258308898" name="Graphic 5" descr="Mobile phone with the words What's Next Text in a speech bubble."><a:extLst><a:ext uri="{C183D7F6-B498-43B3-948B-1728B52AA6E4}"><adec:decorative xmlns:adec="http://schemas.microsoft.com/office/drawing/2017/decorative" val="0"/></a:ext></a:extLst></pic:cNvPr><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId5" cstate="print"><a:extLst><a:ext uri="{BEBA8EAE-BF5A-486C-A8C5-ECC9F3942E4B}"><a14:imgProps xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main"><a14:imgLayer r:embed="rId6"><a14:imgEffect><a14:brightnessContrast bright="-23000"/></a14:imgEffect></a14:imgLayer></a14:imgProps></a:ext><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="987552" cy="1216152"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic><wp14:sizeRelH relativeFrom="margin"><wp14:pctWidth>0</wp14:pctWidth></wp14:sizeRelH><wp14:sizeRelV relativeFrom="margin"><wp14:pctHeight>0</wp14:pctHeight></wp14:sizeRelV></wp:anchor></w:drawing></w:r></w:p>
'@

$r1 = $d.Content
$found1 = $r1.Find.Execute("has_whats_next_text", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not locate the 'has_whats_next_text' anchor text"
}
$ifParaIdx = Get-ParagraphIndexAt $d $r1.Start
if ($ifParaIdx -lt 0) {
    throw "Could not resolve the paragraph containing 'has_whats_next_text'"
}
$picParaIdx = $ifParaIdx + 1
$picPara = $d.Paragraphs.Item($picParaIdx)
$picPara.Range.InsertXML($anchorXml)

# --- Part 2: swap the "{%p endif %}" run content / formatting between the two paragraphs ---

$para1Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="16D02625" w14:textId="4B84AA5D" w:rsidR="00064D95" w:rsidRDefault="00064D95" w:rsidP="00064D95"><w:pPr><w:spacing w:after="0"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">{%p </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>endif %}</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>
'@
$para2Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="4A2AA6A3" w14:textId="439FB7C0" w:rsidR="00064D95" w:rsidRPr="00064D95" w:rsidRDefault="00064D95" w:rsidP="00064D95"><w:pPr><w:spacing w:after="0"/></w:pPr></w:p>
'@

$r2 = $d.Content
$found2 = $r2.Find.Execute("{%p endif %}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not locate the '{%p endif %}' template text"
}
$endifParaIdx = Get-ParagraphIndexAt $d $r2.Start
if ($endifParaIdx -lt 0) {
    throw "Could not resolve the paragraph containing '{%p endif %}'"
}
$blankParaIdx = $endifParaIdx - 1

$blankPara = $d.Paragraphs.Item($blankParaIdx)
$blankPara.Range.InsertXML($para1Xml)

$endifPara = $d.Paragraphs.Item($endifParaIdx)
$endifPara.Range.InsertXML($para2Xml)

Write-Output "edit complete"
